$wb = $excel.ActiveWorkbook

# --- "Repayment schedule" sheet: insert a new (blank) column before column N ---
# This pushes the existing "Late" / second "heading" / "Outstanding" columns
# one position to the right (N->O, O->P, P->Q) and leaves a blank column N,
# matching Excel's normal "Insert Column" behaviour (new column inherits the
# width/style of the column to its left).
$wsSchedule = $wb.Worksheets.Item("Repayment schedule")
$precedingColumnWidth = $wsSchedule.Columns.Item(13).ColumnWidth
$wsSchedule.Columns.Item(14).Insert()
$wsSchedule.Columns.Item(14).ColumnWidth = $precedingColumnWidth

# --- Update the active sheet / selections to match the new workbook state ---
# "Transactions" loses its tabSelected flag and gets a new selection.
$wsTransactions = $wb.Worksheets.Item("Transactions")
$wsTransactions.Activate()
$wsTransactions.Range("D22").Select()

# "Repayment schedule" becomes the active / selected tab with a new selection.
$wsSchedule.Activate()
$wsSchedule.Range("R11").Select()
